# Update countries & provincias Spain
# Applies the data refresh to the "Pais" COVID dashboard sheet:
#  - bumps the "Datos actualizados..." timestamp
#  - refreshes total/active/recovered/critical/deaths counters for a batch of countries
#  - two pairs of countries swapped rank (their row data is exchanged) because
#    one overtook the other in total cases

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    if ($Country -ne $null) {
        $ws.Cells.Item($Row, 1).Value = $Country
    }
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 20:52"

# --- Straightforward numeric refreshes (country stays on the same row) ---
Set-Row 4   $null 6952136 26195 4203484 2545163 0 322  203489  # Estados Unidos
Set-Row 5   $null 5392666 87191 4295755 1010212 0 1074 86699   # India
Set-Row 15  $null 442194  13498 91574   319346  0 26   31274   # Francia
Set-Row 31  $null 125620  1491  97063   17473   0 40   11084   # Ecuador
Set-Row 82  $null 20431   60    19124   891     0 0    416     # Camerun
Set-Row 100 $null 9649    81    8188    1428    0 0    33      # Maldivas
Set-Row 113 $null 5718    2     4030    1509    0 0    179     # Malaui
Set-Row 120 $null 4986    6     3887    1010    0 0    89      # Congo
Set-Row 138 $null 3401    11    2812    491     0 0    98      # Somalia
Set-Row 146 $null 2642    33    1290    1303    0 0    49      # Sudan del Sur

# --- Rank swaps: countries exchange rows, each keeping/getting its own updated stats ---

# Rows 117-118: Cuba overtakes Hong Kong
Set-Row 117 "Cuba"      5055 51 4284 658 0 2 113
Set-Row 118 "Hong Kong" 5010 13 4707 200 0 0 103

# Rows 204-205: Santa Lucia now listed ahead of Timor Oriental (tied totals)
Set-Row 204 "Santa Lucia"    27 0 26 1 0 0 0
Set-Row 205 "Timor Oriental" 27 0 26 1 0 0 0

# Rows 214-215: Montserrat now listed ahead of Islas Malvinas
Set-Row 214 "Montserrat"     13 0 12 0 0 0 1
Set-Row 215 "Islas Malvinas" 13 0 13 0 0 0 0
